# Apply the "first complete version of the manuscript" update to the
# example_data workbook's "Data" sheet: fill in previously-blank
# experiment temperature (T) / pressure (U) columns and correct a few
# existing values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 4: correct temperature reading
$ws.Range("T4").Value = 100

# Rows 5-7: newly measured T/U pair (same conditions)
$ws.Range("T5").Value = 2000
$ws.Range("U5").Value = 1200
$ws.Range("T6").Value = 2000
$ws.Range("U6").Value = 1200
$ws.Range("T7").Value = 2000
$ws.Range("U7").Value = 1200

# Rows 11-13: newly measured T/U pair (same conditions)
$ws.Range("T11").Value = 300
$ws.Range("U11").Value = 900
$ws.Range("T12").Value = 300
$ws.Range("U12").Value = 900
$ws.Range("T13").Value = 300
$ws.Range("U13").Value = 900

# Row 14: only the temperature was recorded
$ws.Range("T14").Value = 300

# Row 15: only the pressure was recorded
$ws.Range("U15").Value = 900

# Rows 16 and 18: corrected temperature values
$ws.Range("T16").Value = 1500
$ws.Range("T18").Value = 1500

# Update the active selection to match the author's saved cursor position
$ws.Activate() | Out-Null
$ws.Range("T19").Select() | Out-Null
